# Auto-generated edit script applying the cryptos.xlsx data refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.349.22'
$ws.Range('E2').Value = '  +0.90%  '
$ws.Range('D3').Value = '1.825.68'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.61'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4473'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3777'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07447'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8901'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.14%  '
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('D12').Value = '1.827.00'
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.747'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.468'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.62'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.07131'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008794'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.15'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.83%  '
$ws.Range('D21').Value = '27.351.51'
$ws.Range('E21').Value = '  +0.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.397'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.73%  '
$ws.Range('E23').Value = '  -0.35%  '
$ws.Range('D24').Value = '2.053.66'
$ws.Range('E24').Value = '  +0.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.970'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.49'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.334'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.82%  '
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.393'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '117.91'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08882'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7956'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.203'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.614'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.96%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.919'
$ws.Range('D35').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.113'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01987'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05316'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.301'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5365'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.870'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.47%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.335'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +17.14%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1718'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.680'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.33%  '
$ws.Range('E46').Value = '  -2.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.61'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.698'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.000'
$ws.Range('D50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06410'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.35%  '
